$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data validation pipeline fix: session for user 2 (vi) on 01/06/2025 was
# re-classified from "Anomala" to "Normal" after review.
$ws.Range("G3").Value = "Normal"

# New sessions ingested from the 06/06/2025 pipeline run.
# Force column D to text first so the dd/mm/yyyy strings aren't
# auto-coerced into date serials (matches the text formatting used by
# the existing rows in this column).
$ws.Range("D10:D12").NumberFormat = "@"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "aline"
$ws.Range("C10").Value = "aline@gmail.com"
$ws.Range("D10").Value = "06/06/2025"
$ws.Range("E10").Value = 70.25
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = "Anômala"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "vih"
$ws.Range("C11").Value = "vih@gmail.com"
$ws.Range("D11").Value = "06/06/2025"
$ws.Range("E11").Value = 8.983333333333333
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = "Normal"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "cad"
$ws.Range("C12").Value = "cad@gmail.com"
$ws.Range("D12").Value = "06/06/2025"
$ws.Range("E12").Value = 0.2833333333333333
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = "Normal"

# Drop the temporary text format again so the new rows carry no explicit
# cell style, matching the rest of the table (header aside).
$ws.Range("D10:D12").ClearFormats()
